$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Karamaji Eliphasi", 5, 18, 1, 3, 19),
    @("Katimbo Dennis Ronald", 6, 27, 0, 7, 11),
    @("Ekiru David Bernard", 4, 21, 0, 15, 13),
    @("Muwanguzi Enock", 1, 14, 0, 0, 0),
    @("Ssali Musa", 5, 15, 0, 2, 4),
    @("Kizza Harper Stephen", 2, 26, 5, 1, 1),
    @("Kwesiga Norman", 0, 1, 1, 1, 0),
    @("Ssekanyonyi Kennedy", 5, 35, 3, 1, 7),
    @("kibalama frank", 11, 45, 4, 2, 1),
    @("Mukasa George Parker", 3, 19, 0, 5, 7),
    @("Barya Mwebaze", 6, 33, 1, 4, 8),
    @("Makumbi Patrick", 7, 22, 3, 2, 2),
    @("Sserunyigo Enock", 15, 31, 2, 1, 5),
    @("David Oluka", 0, 21, 0, 0, 17),
    @("Kiddawalime Jimmy", 0, 1, 0, 0, 0),
    @("Mukasa joseph mutesasira", 3, 6, 0, 0, 0),
    @("Wateya Job", 2, 21, 0, 1, 9),
    @("Musaazi Brian", 7, 2, 0, 3, 3),
    @("Mubiru Stephen", 1, 4, 0, 0, 0),
    @("Kwesiga Samuel", 6, 19, 0, 17, 15),
    @("Aojan James Patrick", 0, 0, 0, 0, 0),
    @("Kyeyune Eric Brian", 0, 10, 0, 1, 5),
    @("NAKIMERA NASHIM", 0, 0, 0, 0, 0),
    @("Musasizi Eric John", 6, 19, 0, 1, 1),
    @("Muyanja Moses", 5, 36, 0, 0, 0),
    @("MBAZIIRA HUDHAIFA", 8, 27, 0, 0, 0),
    @("KIWANDA PETER", 0, 2, 0, 0, 3),
    @("Daphine Nakanwagi", 4, 25, 0, 0, 0),
    @("Lubega Mark", 0, 4, 0, 0, 0),
    @("Talemwa Steven", 0, 11, 0, 0, 0),
    @("other", 0, 0, 0, 0, 0),
    @("Bbaale  Jean", 3, 26, 0, 3, 2),
    @("Boas Massinde", 0, 0, 0, 0, 0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
